# DTT-Test-Hour-Log.xlsx : "Implemented all the feedback. Version 2.0"
# Fill in the previously-blank rows 13-21 on the (only/active) worksheet
# with the new hour-log entries, then move the active-cell selection to A21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13
$ws.Range("A13").Value = "Added the Scrollview"
$ws.Range("B13").Value = 43089
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = "Implenting the scrollview in the about section"

# Row 14 (C14 is logged as text "0.50", not a number, so it is excluded
# from the B30 SUM total further down - match that by round-tripping the
# cell through a text number-format so the apostrophe-less numeric-looking
# string is kept as a literal string instead of being parsed as a number)
$ws.Range("A14").Value = "Research Traits"
$ws.Range("B14").Value = 43089
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "0.50"
$ws.Range("C14").NumberFormat = "0"
$ws.Range("D14").Value = "I didnt know what this was. So i started to learn how it worked."

# Row 15
$ws.Range("A15").Value = "Traits for the aboutview"
$ws.Range("B15").Value = 43089
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "1.50"
$ws.Range("C15").NumberFormat = "0"
$ws.Range("D15").Value = "Traits for every device for the aboutview"

# Row 16
$ws.Range("A16").Value = "Traits for the homeview"
$ws.Range("B16").Value = 43089
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "1.50"
$ws.Range("C16").NumberFormat = "0"
$ws.Range("D16").Value = "Traits for every device for the homeviewcontroller"

# Row 17
$ws.Range("A17").Value = "Traits for the mapview"
$ws.Range("B17").Value = 43090
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "1.50"
$ws.Range("C17").NumberFormat = "0"
$ws.Range("D17").Value = "Traits for every device for the mapviewcontroller"

# Row 18
$ws.Range("A18").Value = "Appstructure"
$ws.Range("B18").Value = 43090
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "0.50"
$ws.Range("C18").NumberFormat = "0"
$ws.Range("D18").Value = "The mapstructure of the app is now better"

# Row 19
$ws.Range("A19").Value = "GPS permissions"
$ws.Range("B19").Value = 43090
$ws.Range("C19").Value = 2
$ws.Range("D19").Value = "Asking with alertviews for permission of the gps"

# Row 20
$ws.Range("A20").Value = "Traits failures"
$ws.Range("B20").Value = 43090
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = "The traits didnt worked that well. I needed to set them again for several times."

# Row 21
$ws.Range("A21").Value = "Internet connection"
$ws.Range("B21").Value = 43090
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "0.5"
$ws.Range("C21").NumberFormat = "0"
$ws.Range("D21").Value = "I searched and tried to understand what i needed to do. But i didnt know where i needed to check."

# Move the selection, matching the recorded sheet view state
$ws.Range("A21").Select()
